$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1285113333333333
$ws.Range("H2").Value = 0.385534
$ws.Range("I2").Value = 0.03749201237720504
$ws.Range("J2").Value = 0.03749201237720504
$ws.Range("M2").Value = 4.441418
$ws.Range("N2").Value = 13.324254
$ws.Range("O2").Value = 0.5990512116078841
$ws.Range("P2").Value = 0.5990512116078841
$ws.Range("Q2").Value = 0.5707725490706667
$ws.Range("R2").Value = 5.136952941635999
$ws.Range("S2").Value = 0.02245963544018247
$ws.Range("T2").Value = 0.02245963544018247
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1285113333333333
$ws.Range("H3").Value = 0.385534
$ws.Range("I3").Value = 0.03749201237720504
$ws.Range("J3").Value = 0.03749201237720504
$ws.Range("M3").Value = 0.9983063333333334
$ws.Range("O3").Value = 0.1346499290404906
$ws.Range("P3").Value = 0.1346499290404906
$ws.Range("Q3").Value = 0.1282936779717778
$ws.Range("R3").Value = 1.154643101746
$ws.Range("S3").Value = 0.005048296806175852
$ws.Range("T3").Value = 0.005048296806175853
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1285113333333333
$ws.Range("H4").Value = 0.385534
$ws.Range("I4").Value = 0.03749201237720504
$ws.Range("J4").Value = 0.03749201237720504
$ws.Range("M4").Value = 0.643907
$ws.Range("N4").Value = 1.931721
$ws.Range("O4").Value = 0.08684912532727113
$ws.Range("P4").Value = 0.08684912532727113
$ws.Range("Q4").Value = 0.08274934711266667
$ws.Range("R4").Value = 0.744744124014
$ws.Range("S4").Value = 0.003256148481719481
$ws.Range("T4").Value = 0.003256148481719481
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1285113333333333
$ws.Range("H5").Value = 0.385534
$ws.Range("I5").Value = 0.03749201237720504
$ws.Range("J5").Value = 0.03749201237720504
$ws.Range("M5").Value = 1.330456
$ws.Range("N5").Value = 3.991368
$ws.Range("O5").Value = 0.1794497340243542
$ws.Range("P5").Value = 0.1794497340243542
$ws.Range("Q5").Value = 0.1709786745013334
$ws.Range("R5").Value = 1.538808070512
$ws.Range("S5").Value = 0.00672793164912724
$ws.Range("T5").Value = 0.00672793164912724
$ws.Range("I6").Value = 0.7552862722193517
$ws.Range("J6").Value = 0.755286272219352
$ws.Range("M6").Value = 4.441418
$ws.Range("N6").Value = 13.324254
$ws.Range("O6").Value = 0.5990512116078841
$ws.Range("P6").Value = 0.5990512116078841
$ws.Range("Q6").Value = 11.498360411692
$ws.Range("R6").Value = 103.485243705228
$ws.Range("S6").Value = 0.4524551564838049
$ws.Range("T6").Value = 0.452455156483805
$ws.Range("I7").Value = 0.7552862722193517
$ws.Range("J7").Value = 0.755286272219352
$ws.Range("M7").Value = 0.9983063333333334
$ws.Range("O7").Value = 0.1346499290404906
$ws.Range("P7").Value = 0.1346499290404906
$ws.Range("Q7").Value = 2.584509276528667
$ws.Range("S7").Value = 0.1016992429595924
$ws.Range("T7").Value = 0.1016992429595924
$ws.Range("I8").Value = 0.7552862722193517
$ws.Range("J8").Value = 0.755286272219352
$ws.Range("M8").Value = 0.643907
$ws.Range("N8").Value = 1.931721
$ws.Range("O8").Value = 0.08684912532727113
$ws.Range("P8").Value = 0.08684912532727113
$ws.Range("Q8").Value = 1.667006968858
$ws.Range("R8").Value = 15.003062719722
$ws.Range("S8").Value = 0.06559595211394589
$ws.Range("T8").Value = 0.06559595211394592
$ws.Range("I9").Value = 0.7552862722193517
$ws.Range("J9").Value = 0.755286272219352
$ws.Range("M9").Value = 1.330456
$ws.Range("N9").Value = 3.991368
$ws.Range("O9").Value = 0.1794497340243542
$ws.Range("P9").Value = 0.1794497340243542
$ws.Range("Q9").Value = 3.444409555664
$ws.Range("R9").Value = 30.999686000976
$ws.Range("S9").Value = 0.1355359206620087
$ws.Range("T9").Value = 0.1355359206620087
$ws.Range("G10").Value = 0.692415
$ws.Range("H10").Value = 2.077245
$ws.Range("I10").Value = 0.2020057770533527
$ws.Range("J10").Value = 0.2020057770533527
$ws.Range("M10").Value = 4.441418
$ws.Range("N10").Value = 13.324254
$ws.Range("O10").Value = 0.5990512116078841
$ws.Range("P10").Value = 0.5990512116078841
$ws.Range("Q10").Value = 3.07530444447
$ws.Range("R10").Value = 27.67774000023
$ws.Range("S10").Value = 0.1210118054956031
$ws.Range("T10").Value = 0.1210118054956031
$ws.Range("G11").Value = 0.692415
$ws.Range("H11").Value = 2.077245
$ws.Range("I11").Value = 0.2020057770533527
$ws.Range("J11").Value = 0.2020057770533527
$ws.Range("M11").Value = 0.9983063333333334
$ws.Range("O11").Value = 0.1346499290404906
$ws.Range("P11").Value = 0.1346499290404906
$ws.Range("Q11").Value = 0.6912422797950001
$ws.Range("R11").Value = 6.221180518155001
$ws.Range("S11").Value = 0.0272000635460031
$ws.Range("T11").Value = 0.02720006354600311
$ws.Range("G12").Value = 0.692415
$ws.Range("H12").Value = 2.077245
$ws.Range("I12").Value = 0.2020057770533527
$ws.Range("J12").Value = 0.2020057770533527
$ws.Range("M12").Value = 0.643907
$ws.Range("N12").Value = 1.931721
$ws.Range("O12").Value = 0.08684912532727113
$ws.Range("P12").Value = 0.08684912532727113
$ws.Range("Q12").Value = 0.445850865405
$ws.Range("R12").Value = 4.012657788645
$ws.Range("S12").Value = 0.01754402504813942
$ws.Range("T12").Value = 0.01754402504813942
$ws.Range("G13").Value = 0.692415
$ws.Range("H13").Value = 2.077245
$ws.Range("I13").Value = 0.2020057770533527
$ws.Range("J13").Value = 0.2020057770533527
$ws.Range("M13").Value = 1.330456
$ws.Range("N13").Value = 3.991368
$ws.Range("O13").Value = 0.1794497340243542
$ws.Range("P13").Value = 0.1794497340243542
$ws.Range("Q13").Value = 0.92122769124
$ws.Range("R13").Value = 8.291049221160002
$ws.Range("S13").Value = 0.03624988296360714
$ws.Range("T13").Value = 0.03624988296360714
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.01787866666666667
$ws.Range("H14").Value = 0.053636
$ws.Range("I14").Value = 0.005215938350090445
$ws.Range("J14").Value = 0.005215938350090446
$ws.Range("M14").Value = 4.441418
$ws.Range("N14").Value = 13.324254
$ws.Range("O14").Value = 0.5990512116078841
$ws.Range("P14").Value = 0.5990512116078841
$ws.Range("Q14").Value = 0.07940663194933333
$ws.Range("R14").Value = 0.714659687544
$ws.Range("S14").Value = 0.003124614188293709
$ws.Range("T14").Value = 0.00312461418829371
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.01787866666666667
$ws.Range("H15").Value = 0.053636
$ws.Range("I15").Value = 0.005215938350090445
$ws.Range("J15").Value = 0.005215938350090446
$ws.Range("M15").Value = 0.9983063333333334
$ws.Range("O15").Value = 0.1346499290404906
$ws.Range("P15").Value = 0.1346499290404906
$ws.Range("Q15").Value = 0.01784838616488889
$ws.Range("R15").Value = 0.160635475484
$ws.Range("S15").Value = 0.000702325728719252
$ws.Range("T15").Value = 0.0007023257287192521
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.01787866666666667
$ws.Range("H16").Value = 0.053636
$ws.Range("I16").Value = 0.005215938350090445
$ws.Range("J16").Value = 0.005215938350090446
$ws.Range("M16").Value = 0.643907
$ws.Range("N16").Value = 1.931721
$ws.Range("O16").Value = 0.08684912532727113
$ws.Range("P16").Value = 0.08684912532727113
$ws.Range("Q16").Value = 0.01151219861733333
$ws.Range("R16").Value = 0.103609787556
$ws.Range("S16").Value = 0.0004529996834663249
$ws.Range("T16").Value = 0.0004529996834663249
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.01787866666666667
$ws.Range("H17").Value = 0.053636
$ws.Range("I17").Value = 0.005215938350090445
$ws.Range("J17").Value = 0.005215938350090446
$ws.Range("M17").Value = 1.330456
$ws.Range("N17").Value = 3.991368
$ws.Range("O17").Value = 0.1794497340243542
$ws.Range("P17").Value = 0.1794497340243542
$ws.Range("Q17").Value = 0.02378677933866667
$ws.Range("R17").Value = 0.214081014048
$ws.Range("S17").Value = 0.0009359987496111593
$ws.Range("T17").Value = 0.0009359987496111593
